# Update the Il34-Ptprz1 LR-pairs export to the refreshed TPM numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" sending-cluster block (old rows 17-19) is dropped entirely.
$ws.Range("A17:T19").Delete()

# The old "Neutrophils" block (rows 14-16) becomes "Resolving-Mac" in the refreshed export;
# this also drops the now-unused "Neutrophils" shared string from the workbook.
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("A16").Value = "Resolving-Mac"

# Refresh the ligand/receptor expression metrics (columns E:T) for every remaining row.
# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.957220666666667
$ws.Range("H2").Value = 14.871662
$ws.Range("I2").Value = 0.3393483464618172
$ws.Range("J2").Value = 0.3393483464618172
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04936366666666667
$ws.Range("N2").Value = 0.148091
$ws.Range("O2").Value = 0.04616170608573571
$ws.Range("P2").Value = 0.0461617060857357
$ws.Range("Q2").Value = 0.2447065885824445
$ws.Range("R2").Value = 2.202359297242
$ws.Range("S2").Value = 0.01566489863005082
$ws.Range("T2").Value = 0.01566489863005081

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.957220666666667
$ws.Range("H3").Value = 14.871662
$ws.Range("I3").Value = 0.3393483464618172
$ws.Range("J3").Value = 0.3393483464618172
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.008616666666666667
$ws.Range("N3").Value = 0.02585
$ws.Range("O3").Value = 0.008057748967298944
$ws.Range("P3").Value = 0.008057748967298944
$ws.Range("Q3").Value = 0.04271471807777778
$ws.Range("R3").Value = 0.3844324627
$ws.Range("S3").Value = 0.002734383788257312
$ws.Range("T3").Value = 0.002734383788257312

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.957220666666667
$ws.Range("H4").Value = 14.871662
$ws.Range("I4").Value = 0.3393483464618172
$ws.Range("J4").Value = 0.3393483464618172
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.011383666666666
$ws.Range("N4").Value = 3.034151
$ws.Range("O4").Value = 0.9457805449469654
$ws.Range("P4").Value = 0.9457805449469653
$ws.Range("Q4").Value = 5.013652014329111
$ws.Range("R4").Value = 45.122868128962
$ws.Range("S4").Value = 0.3209490640435091
$ws.Range("T4").Value = 0.3209490640435091

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.165549333333333
$ws.Range("H5").Value = 6.496648
$ws.Range("I5").Value = 0.1482434684398066
$ws.Range("J5").Value = 0.1482434684398067
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04936366666666667
$ws.Range("N5").Value = 0.148091
$ws.Range("O5").Value = 0.04616170608573571
$ws.Range("P5").Value = 0.0461617060857357
$ws.Range("Q5").Value = 0.1068994554408889
$ws.Range("R5").Value = 0.9620950989680001
$ws.Range("S5").Value = 0.006843171419248391
$ws.Range("T5").Value = 0.006843171419248391

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.165549333333333
$ws.Range("H6").Value = 6.496648
$ws.Range("I6").Value = 0.1482434684398066
$ws.Range("J6").Value = 0.1482434684398067
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.008616666666666667
$ws.Range("N6").Value = 0.02585
$ws.Range("O6").Value = 0.008057748967298944
$ws.Range("P6").Value = 0.008057748967298944
$ws.Range("Q6").Value = 0.01865981675555555
$ws.Range("R6").Value = 0.1679383508
$ws.Range("S6").Value = 0.001194508654729665
$ws.Range("T6").Value = 0.001194508654729666

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.165549333333333
$ws.Range("H7").Value = 6.496648
$ws.Range("I7").Value = 0.1482434684398066
$ws.Range("J7").Value = 0.1482434684398067
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.011383666666666
$ws.Range("N7").Value = 3.034151
$ws.Range("O7").Value = 0.9457805449469654
$ws.Range("P7").Value = 0.9457805449469653
$ws.Range("Q7").Value = 2.190201225094222
$ws.Range("R7").Value = 19.711811025848
$ws.Range("S7").Value = 0.1402057883658286
$ws.Range("T7").Value = 0.1402057883658286

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.03111233333333334
$ws.Range("H8").Value = 0.093337
$ws.Range("I8").Value = 0.002129806111361772
$ws.Range("J8").Value = 0.002129806111361772
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04936366666666667
$ws.Range("N8").Value = 0.148091
$ws.Range("O8").Value = 0.04616170608573571
$ws.Range("P8").Value = 0.0461617060857357
$ws.Range("Q8").Value = 0.001535818851888889
$ws.Range("R8").Value = 0.013822369667
$ws.Range("S8").Value = [double]"9.83154837322858E-05"
$ws.Range("T8").Value = [double]"9.831548373228579E-05"

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.03111233333333334
$ws.Range("H9").Value = 0.093337
$ws.Range("I9").Value = 0.002129806111361772
$ws.Range("J9").Value = 0.002129806111361772
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.008616666666666667
$ws.Range("N9").Value = 0.02585
$ws.Range("O9").Value = 0.008057748967298944
$ws.Range("P9").Value = 0.008057748967298944
$ws.Range("Q9").Value = 0.0002680846055555555
$ws.Range("R9").Value = 0.00241276145
$ws.Range("S9").Value = [double]"1.71614429943723E-05"
$ws.Range("T9").Value = [double]"1.71614429943723E-05"

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.03111233333333334
$ws.Range("H10").Value = 0.093337
$ws.Range("I10").Value = 0.002129806111361772
$ws.Range("J10").Value = 0.002129806111361772
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.011383666666666
$ws.Range("N10").Value = 3.034151
$ws.Range("O10").Value = 0.9457805449469654
$ws.Range("P10").Value = 0.9457805449469653
$ws.Range("Q10").Value = 0.03146650576522222
$ws.Range("R10").Value = 0.283198551887
$ws.Range("S10").Value = 0.002014329184635114
$ws.Range("T10").Value = 0.002014329184635114

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 7.327472333333333
$ws.Range("H11").Value = 21.982417
$ws.Range("I11").Value = 0.5016047876951574
$ws.Range("J11").Value = 0.5016047876951574
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04936366666666667
$ws.Range("N11").Value = 0.148091
$ws.Range("O11").Value = 0.04616170608573571
$ws.Range("P11").Value = 0.0461617060857357
$ws.Range("Q11").Value = 0.3617109017718889
$ws.Range("R11").Value = 3.255398115947
$ws.Range("S11").Value = 0.02315493278078171
$ws.Range("T11").Value = 0.02315493278078171

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 7.327472333333333
$ws.Range("H12").Value = 21.982417
$ws.Range("I12").Value = 0.5016047876951574
$ws.Range("J12").Value = 0.5016047876951574
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.008616666666666667
$ws.Range("N12").Value = 0.02585
$ws.Range("O12").Value = 0.008057748967298944
$ws.Range("P12").Value = 0.008057748967298944
$ws.Range("Q12").Value = 0.06313838660555555
$ws.Range("R12").Value = 0.56824547945
$ws.Range("S12").Value = 0.00404180546004286
$ws.Range("T12").Value = 0.00404180546004286

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 7.327472333333333
$ws.Range("H13").Value = 21.982417
$ws.Range("I13").Value = 0.5016047876951574
$ws.Range("J13").Value = 0.5016047876951574
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.011383666666666
$ws.Range("N13").Value = 3.034151
$ws.Range("O13").Value = 0.9457805449469654
$ws.Range("P13").Value = 0.9457805449469653
$ws.Range("Q13").Value = 7.41088583588522
$ws.Range("R13").Value = 66.69797252296699
$ws.Range("S13").Value = 0.4744080494543328
$ws.Range("T13").Value = 0.4744080494543328

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.1267043333333333
$ws.Range("H14").Value = 0.380113
$ws.Range("I14").Value = 0.008673591291857003
$ws.Range("J14").Value = 0.008673591291857003
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.04936366666666667
$ws.Range("N14").Value = 0.148091
$ws.Range("O14").Value = 0.04616170608573571
$ws.Range("P14").Value = 0.0461617060857357
$ws.Range("Q14").Value = 0.006254590475888889
$ws.Range("R14").Value = 0.056291314283
$ws.Range("S14").Value = 0.0004003877719224996
$ws.Range("T14").Value = 0.0004003877719224996

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.1267043333333333
$ws.Range("H15").Value = 0.380113
$ws.Range("I15").Value = 0.008673591291857003
$ws.Range("J15").Value = 0.008673591291857003
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.008616666666666667
$ws.Range("N15").Value = 0.02585
$ws.Range("O15").Value = 0.008057748967298944
$ws.Range("P15").Value = 0.008057748967298944
$ws.Range("Q15").Value = 0.001091769005555556
$ws.Range("R15").Value = 0.00982592105
$ws.Range("S15").Value = [double]"6.988962127473388E-05"
$ws.Range("T15").Value = [double]"6.988962127473388E-05"

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.1267043333333333
$ws.Range("H16").Value = 0.380113
$ws.Range("I16").Value = 0.008673591291857003
$ws.Range("J16").Value = 0.008673591291857003
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.011383666666666
$ws.Range("N16").Value = 3.034151
$ws.Range("O16").Value = 0.9457805449469654
$ws.Range("P16").Value = 0.9457805449469653
$ws.Range("Q16").Value = 0.1281466932292222
$ws.Range("R16").Value = 1.153320239063
$ws.Range("S16").Value = 0.00820331389865977
$ws.Range("T16").Value = 0.00820331389865977

